# Adds the submitter's data row (name, email, repo link) to the
# "Open Source task" workbook, mirroring the author's commit
# "I have added my data":
#   A2 -> name (Arabic)
#   B2 -> email, turned into a mailto: hyperlink (Excel's usual
#         auto-hyperlink behaviour when an e-mail address is typed in)
#   C2 -> the GitHub repo link (plain text, as in the source diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$name  = "محمد متولي ابوالنصر متولي سالم"
$email = "mohamedhigazii60@gmail.com"
$repo  = "https://github.com/mohamedhegazy03/Open-Source-.git"

$ws.Range("A2").Value = $name
$ws.Range("B2").Value = $email
$ws.Range("C2").Value = $repo

# Turn the e-mail cell into a hyperlink (adds the Hyperlink style/font
# and the relationship automatically, just like Excel does when you
# insert a hyperlink over existing text).
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:" + $email)

# The new values are wider than the original columns, so let Excel
# resize ("best fit") the columns to fit the new content.
$ws.Columns("A:C").AutoFit() | Out-Null

# Put the selection where the user's cursor ended up after entering
# the new row of data.
$ws.Range("C13").Select() | Out-Null
